$wb = $excel.ActiveWorkbook

# --- Mix_Mass_Frac sheet ---
$ws = $wb.Worksheets.Item("Mix_Mass_Frac")
$ws.Rows.Item(11).Insert()
$ws.Cells.Item(2, 1).Value = 2
$ws.Cells.Item(2, 2).Value = 0.7165025054715928
$ws.Cells.Item(2, 3).Value = 0.06863842046557493
$ws.Cells.Item(2, 4).Value = 0.0007990841499006221
$ws.Cells.Item(2, 5).Value = 0.03430951048139613
$ws.Cells.Item(2, 6).Value = 0.0030834031642105702
$ws.Cells.Item(2, 7).Value = 0.00032962744917057626
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 0.1932701325617664
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 0.011279529457069951
$ws.Cells.Item(2, 13).Value = 0.634242087913461
$ws.Cells.Item(2, 14).Value = 0.03641910327615504
$ws.Cells.Item(2, 15).Value = 0.011914908462461746
$ws.Cells.Item(2, 16).Value = 0
$ws.Cells.Item(2, 17).Value = 0.0009375491827481847
$ws.Cells.Item(2, 18).Value = 0.0018747004808270912
$ws.Cells.Item(2, 19).Value = 0.002901942955257797
$ws.Cells.Item(2, 20).Value = 0.0006075382495578973
$ws.Cells.Item(3, 1).Value = 4
$ws.Cells.Item(3, 2).Value = 0.6510128786562844
$ws.Cells.Item(3, 3).Value = 0.08395163995473542
$ws.Cells.Item(3, 4).Value = 0.000997639124510761
$ws.Cells.Item(3, 5).Value = 0.03859024236548231
$ws.Cells.Item(3, 6).Value = 0.0038475542626794017
$ws.Cells.Item(3, 7).Value = 0.00008776046843684441
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 0.24229601012811072
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 0.014395326538266374
$ws.Cells.Item(3, 13).Value = 0.546946552259967
$ws.Cells.Item(3, 14).Value = 0.04670750803171299
$ws.Cells.Item(3, 15).Value = 0.014503468454741117
$ws.Cells.Item(3, 16).Value = 0
$ws.Cells.Item(3, 17).Value = 0.0014136335505368191
$ws.Cells.Item(3, 18).Value = 0.002436152842924903
$ws.Cells.Item(3, 19).Value = 0.0038265120178952805
$ws.Cells.Item(3, 20).Value = 0.0013244690126094799
$ws.Cells.Item(4, 1).Value = 6
$ws.Cells.Item(4, 2).Value = 0.5475288334865248
$ws.Cells.Item(4, 3).Value = 0.10693683020356214
$ws.Cells.Item(4, 4).Value = 0.0012956404020672083
$ws.Cells.Item(4, 5).Value = 0.04372830188322384
$ws.Cells.Item(4, 6).Value = 0.0052656210920229134
$ws.Cells.Item(4, 7).Value = 0.0002446366016860045
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0.32347423487878285
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0.019013050847986712
$ws.Cells.Item(4, 13).Value = 0.4136482160300555
$ws.Cells.Item(4, 14).Value = 0.058311481814932696
$ws.Cells.Item(4, 15).Value = 0.018546154195211547
$ws.Cells.Item(4, 16).Value = 0
$ws.Cells.Item(4, 17).Value = 0.0019482838823551576
$ws.Cells.Item(4, 18).Value = 0.002805349111804852
$ws.Cells.Item(4, 19).Value = 0.004782199056308547
$ws.Cells.Item(4, 20).Value = 0.0026109593086994154
$ws.Cells.Item(5, 1).Value = 10
$ws.Cells.Item(5, 2).Value = 0.41431406819112365
$ws.Cells.Item(5, 3).Value = 0.1283185967356951
$ws.Cells.Item(5, 4).Value = 0.0015081516002418826
$ws.Cells.Item(5, 5).Value = 0.059489121816044134
$ws.Cells.Item(5, 6).Value = 0.0069522600485077655
$ws.Cells.Item(5, 7).Value = 0.0019355192791594713
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0.42328813296804463
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 0.02218886819733951
$ws.Cells.Item(5, 13).Value = 0.2526943652241502
$ws.Cells.Item(5, 14).Value = 0.07211205470895946
$ws.Cells.Item(5, 15).Value = 0.01997932941222638
$ws.Cells.Item(5, 16).Value = 0
$ws.Cells.Item(5, 17).Value = 0.00249042543049279
$ws.Cells.Item(5, 18).Value = 0.0032052431773109567
$ws.Cells.Item(5, 19).Value = 0.005837931401827756
$ws.Cells.Item(5, 20).Value = 0.0043327587286179905
$ws.Cells.Item(6, 1).Value = 14
$ws.Cells.Item(6, 2).Value = 0.2722270046159259
$ws.Cells.Item(6, 3).Value = 0.13634068049614936
$ws.Cells.Item(6, 4).Value = 0.001632433023356703
$ws.Cells.Item(6, 5).Value = 0.06735912153842395
$ws.Cells.Item(6, 6).Value = 0.009067800251142816
$ws.Cells.Item(6, 7).Value = 0.013941672521544143
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = 0.5424719921817093
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 12).Value = 0.020125956592672013
$ws.Cells.Item(6, 13).Value = 0.11293835281919416
$ws.Cells.Item(6, 14).Value = 0.07243638524071361
$ws.Cells.Item(6, 15).Value = 0.014438892038429688
$ws.Cells.Item(6, 16).Value = 0
$ws.Cells.Item(6, 17).Value = 0.0018737120062576134
$ws.Cells.Item(6, 18).Value = 0.0024933500182573254
$ws.Cells.Item(6, 19).Value = 0.004879651272149198
$ws.Cells.Item(6, 20).Value = 0.0054479141853711185
$ws.Cells.Item(7, 1).Value = 20
$ws.Cells.Item(7, 2).Value = 0.20938448719332356
$ws.Cells.Item(7, 3).Value = 0.13021727854361778
$ws.Cells.Item(7, 4).Value = 0.0013562041710173218
$ws.Cells.Item(7, 5).Value = 0.0610240424977761
$ws.Cells.Item(7, 6).Value = 0.010282461397742729
$ws.Cells.Item(7, 7).Value = 0.02733571217882414
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 0.603533508279936
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 0.017998475069833833
$ws.Cells.Item(7, 13).Value = 0.06647731589801455
$ws.Cells.Item(7, 14).Value = 0.06254028698564834
$ws.Cells.Item(7, 15).Value = 0.010824665781917582
$ws.Cells.Item(7, 16).Value = 0
$ws.Cells.Item(7, 17).Value = 0.0016455035889346744
$ws.Cells.Item(7, 18).Value = 0.002304829188686161
$ws.Cells.Item(7, 19).Value = 0.004459716418050821
$ws.Cells.Item(7, 20).Value = 0.004749909896836446
$ws.Cells.Item(8, 1).Value = 30
$ws.Cells.Item(8, 2).Value = 0.10727721202339385
$ws.Cells.Item(8, 3).Value = 0.09551539168112479
$ws.Cells.Item(8, 4).Value = 0.0008065824289021387
$ws.Cells.Item(8, 5).Value = 0.049336862137483894
$ws.Cells.Item(8, 6).Value = 0.011035600912545213
$ws.Cells.Item(8, 7).Value = 0.07998702598535967
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 10).Value = 0.691253978879161
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 12).Value = 0.00867482611100993
$ws.Cells.Item(8, 13).Value = 0.015372523488123645
$ws.Cells.Item(8, 14).Value = 0.04085963984588869
$ws.Cells.Item(8, 15).Value = 0.0035325665617994963
$ws.Cells.Item(8, 16).Value = 0
$ws.Cells.Item(8, 17).Value = 0.0006265463427197508
$ws.Cells.Item(8, 18).Value = 0.0008830990745864223
$ws.Cells.Item(8, 19).Value = 0.0021153565512951854
$ws.Cells.Item(8, 20).Value = 0.003949034051204005
$ws.Cells.Item(9, 1).Value = 45
$ws.Cells.Item(9, 2).Value = 0.056149852841862254
$ws.Cells.Item(9, 3).Value = 0.07179667711319962
$ws.Cells.Item(9, 4).Value = 0.0002128956555545532
$ws.Cells.Item(9, 5).Value = 0.042039180019614694
$ws.Cells.Item(9, 6).Value = 0.011594908268795131
$ws.Cells.Item(9, 7).Value = 0.11818525920305735
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 10).Value = 0.7268745268424366
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 0.0029552362372269663
$ws.Cells.Item(9, 13).Value = 0.001890074752775993
$ws.Cells.Item(9, 14).Value = 0.022032642373279888
$ws.Cells.Item(9, 15).Value = 0.0013149801784385264
$ws.Cells.Item(9, 16).Value = 0
$ws.Cells.Item(9, 17).Value = 0.000220525098798073
$ws.Cells.Item(9, 18).Value = 0.00021998597260094153
$ws.Cells.Item(9, 19).Value = 0.0006631082842217381
$ws.Cells.Item(9, 20).Value = 0.0013954514318169527
$ws.Cells.Item(10, 1).Value = 60
$ws.Cells.Item(10, 2).Value = 0.030493777848895353
$ws.Cells.Item(10, 3).Value = 0.053270244601418326
$ws.Cells.Item(10, 4).Value = 0.00012379244280364573
$ws.Cells.Item(10, 5).Value = 0.024259918294598743
$ws.Cells.Item(10, 6).Value = 0.011944736397900512
$ws.Cells.Item(10, 7).Value = 0.1584904666308264
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = 0.7431359507846281
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 0.0006700740735281542
$ws.Cells.Item(10, 13).Value = 0.00012299327245527121
$ws.Cells.Item(10, 14).Value = 0.007491967111621844
$ws.Cells.Item(10, 15).Value = 0.0003654231351161034
$ws.Cells.Item(10, 16).Value = 0
$ws.Cells.Item(10, 17).Value = 0
$ws.Cells.Item(10, 18).Value = 0.000022287971743192615
$ws.Cells.Item(10, 19).Value = 0.00010214528335972133
$ws.Cells.Item(10, 20).Value = 0.0003682574030173549
$ws.Cells.Item(11, 1).Value = 75
$ws.Cells.Item(11, 2).Value = 0.017380412021122782
$ws.Cells.Item(11, 3).Value = 0.036681189601409606
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = 0.016684162082519993
$ws.Cells.Item(11, 6).Value = 0.012141612930736128
$ws.Cells.Item(11, 7).Value = 0.183389754119925
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 10).Value = 0.7497924682413492
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 0.00006821602408320213
$ws.Cells.Item(11, 13).Value = 0
$ws.Cells.Item(11, 14).Value = 0.0012425969999767966
$ws.Cells.Item(11, 15).Value = 0
$ws.Cells.Item(11, 16).Value = 0
$ws.Cells.Item(11, 17).Value = 0
$ws.Cells.Item(11, 18).Value = 0
$ws.Cells.Item(11, 19).Value = 0
$ws.Cells.Item(11, 20).Value = 0.00018807664256232414
$ws.Cells.Item(12, 1).Value = 100
$ws.Cells.Item(12, 2).Value = 0.008113561767794101
$ws.Cells.Item(12, 3).Value = 0.018068094369527492
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = 0.007639313586239421
$ws.Cells.Item(12, 6).Value = 0.01219693286797774
$ws.Cells.Item(12, 7).Value = 0.2029596179489102
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 0.7589638693223241
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 0.000011987868108138962
$ws.Cells.Item(12, 13).Value = 0
$ws.Cells.Item(12, 14).Value = 0.00016018403691292216
$ws.Cells.Item(12, 15).Value = 0
$ws.Cells.Item(12, 16).Value = 0
$ws.Cells.Item(12, 17).Value = 0
$ws.Cells.Item(12, 18).Value = 0
$ws.Cells.Item(12, 19).Value = 0
$ws.Cells.Item(12, 20).Value = 0.000024952736895236028

# --- Uncertainties sheet ---
$ws = $wb.Worksheets.Item("Uncertainties")
$ws.Rows.Item(11).Insert()
$ws.Cells.Item(2, 1).Value = 2
$ws.Cells.Item(2, 2).Value = 0.07783603752069966
$ws.Cells.Item(2, 3).Value = 0.0052146980384027465
$ws.Cells.Item(2, 4).Value = 0.00008631770729021734
$ws.Cells.Item(2, 5).Value = 0.007073178094888569
$ws.Cells.Item(2, 6).Value = 0.00034586985105102024
$ws.Cells.Item(2, 7).Value = 0.0004294189372100823
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 0.01568358785728334
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 0.0013177097479420098
$ws.Cells.Item(2, 13).Value = 0.07748476257197127
$ws.Cells.Item(2, 14).Value = 0.00932307455848934
$ws.Cells.Item(2, 15).Value = 0.003014428726877323
$ws.Cells.Item(2, 16).Value = 0
$ws.Cells.Item(2, 17).Value = 0.00021715221609071216
$ws.Cells.Item(2, 18).Value = 0.0002356182119527143
$ws.Cells.Item(2, 19).Value = 0.0002727684193560105
$ws.Cells.Item(2, 20).Value = 0.00024828103977743394
$ws.Cells.Item(3, 1).Value = 4
$ws.Cells.Item(3, 2).Value = 0.05279842851661869
$ws.Cells.Item(3, 3).Value = 0.004627390166709996
$ws.Cells.Item(3, 4).Value = 0.00011983615993016051
$ws.Cells.Item(3, 5).Value = 0.008874340394517318
$ws.Cells.Item(3, 6).Value = 0.0003601907573378255
$ws.Cells.Item(3, 7).Value = 0.000033118931432930436
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 0.012435200387784342
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 0.0010305937178930252
$ws.Cells.Item(3, 13).Value = 0.05222011433322526
$ws.Cells.Item(3, 14).Value = 0.010024045354051879
$ws.Cells.Item(3, 15).Value = 0.003462696453111339
$ws.Cells.Item(3, 16).Value = 0
$ws.Cells.Item(3, 17).Value = 0.00033348985268507133
$ws.Cells.Item(3, 18).Value = 0.00015653114338290557
$ws.Cells.Item(3, 19).Value = 0.0002945608098902118
$ws.Cells.Item(3, 20).Value = 0.00036890985645398694
$ws.Cells.Item(4, 1).Value = 6
$ws.Cells.Item(4, 2).Value = 0.14975078340104936
$ws.Cells.Item(4, 3).Value = 0.013325983773995754
$ws.Cells.Item(4, 4).Value = 0.00016583596212898344
$ws.Cells.Item(4, 5).Value = 0.01054232993935655
$ws.Cells.Item(4, 6).Value = 0.0008958678191639946
$ws.Cells.Item(4, 7).Value = 0.0003880683791059169
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0.04801682498308119
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0.0026527026514707266
$ws.Cells.Item(4, 13).Value = 0.1494736348961462
$ws.Cells.Item(4, 14).Value = 0.008346237415945085
$ws.Cells.Item(4, 15).Value = 0.002919158328210325
$ws.Cells.Item(4, 16).Value = 0
$ws.Cells.Item(4, 17).Value = 0.00046736991630810744
$ws.Cells.Item(4, 18).Value = 0.00042025348462503264
$ws.Cells.Item(4, 19).Value = 0.0006231591641851966
$ws.Cells.Item(4, 20).Value = 0.00014354458945682385
$ws.Cells.Item(5, 1).Value = 10
$ws.Cells.Item(5, 2).Value = 0.08736801799766904
$ws.Cells.Item(5, 3).Value = 0.009524892504705427
$ws.Cells.Item(5, 4).Value = 0.00021028437084416729
$ws.Cells.Item(5, 5).Value = 0.012036061428864507
$ws.Cells.Item(5, 6).Value = 0.00069748007318369
$ws.Cells.Item(5, 7).Value = 0.00016361981513756763
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0.03522487355037707
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 0.0019539468368948714
$ws.Cells.Item(5, 13).Value = 0.08714042253226487
$ws.Cells.Item(5, 14).Value = 0.005437562369031022
$ws.Cells.Item(5, 15).Value = 0.0015293555969379655
$ws.Cells.Item(5, 16).Value = 0
$ws.Cells.Item(5, 17).Value = 0.0005267637642193066
$ws.Cells.Item(5, 18).Value = 0.00026610259548824524
$ws.Cells.Item(5, 19).Value = 0.000601073227261808
$ws.Cells.Item(5, 20).Value = 0.0006087309186854671
$ws.Cells.Item(6, 1).Value = 14
$ws.Cells.Item(6, 2).Value = 0.021780701454469423
$ws.Cells.Item(6, 3).Value = 0.0074310775376446175
$ws.Cells.Item(6, 4).Value = 0.00012122991374817747
$ws.Cells.Item(6, 5).Value = 0.02283326712198656
$ws.Cells.Item(6, 6).Value = 0.0010751258427360326
$ws.Cells.Item(6, 7).Value = 0.013090046750940464
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = 0.04906524998422683
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 12).Value = 0.002499570743260241
$ws.Cells.Item(6, 13).Value = 0.02028474479401484
$ws.Cells.Item(6, 14).Value = 0.008475349108119173
$ws.Cells.Item(6, 15).Value = 0.0035870124288600474
$ws.Cells.Item(6, 16).Value = 0
$ws.Cells.Item(6, 17).Value = 0.00040373460525062164
$ws.Cells.Item(6, 18).Value = 0.0003013623414815196
$ws.Cells.Item(6, 19).Value = 0.0006401656577867455
$ws.Cells.Item(6, 20).Value = 0.0007519470878761389
$ws.Cells.Item(7, 1).Value = 20
$ws.Cells.Item(7, 2).Value = 0.03534301679636145
$ws.Cells.Item(7, 3).Value = 0.011031979069653749
$ws.Cells.Item(7, 4).Value = 0.00021151330261691347
$ws.Cells.Item(7, 5).Value = 0.011885026850516525
$ws.Cells.Item(7, 6).Value = 0.0019987505585005562
$ws.Cells.Item(7, 7).Value = 0.02965779284035285
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 0.08763869588610335
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 0.005461963529707279
$ws.Cells.Item(7, 13).Value = 0.030389678000023932
$ws.Cells.Item(7, 14).Value = 0.02229242141205555
$ws.Cells.Item(7, 15).Value = 0.006055319617956613
$ws.Cells.Item(7, 16).Value = 0
$ws.Cells.Item(7, 17).Value = 0.0008614458701125409
$ws.Cells.Item(7, 18).Value = 0.0005453481176802482
$ws.Cells.Item(7, 19).Value = 0.0012102409831355563
$ws.Cells.Item(7, 20).Value = 0.0001237910839600038
$ws.Cells.Item(8, 1).Value = 30
$ws.Cells.Item(8, 2).Value = 0.020316303454774743
$ws.Cells.Item(8, 3).Value = 0.022239542079306552
$ws.Cells.Item(8, 4).Value = 0.00032318456715713285
$ws.Cells.Item(8, 5).Value = 0.01708392094461895
$ws.Cells.Item(8, 6).Value = 0.0011530575181502485
$ws.Cells.Item(8, 7).Value = 0.036953147451188
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 10).Value = 0.058650807457397126
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 12).Value = 0.003174997297723455
$ws.Cells.Item(8, 13).Value = 0.003415568501192774
$ws.Cells.Item(8, 14).Value = 0.02457525135632295
$ws.Cells.Item(8, 15).Value = 0.0015794596747419384
$ws.Cells.Item(8, 16).Value = 0
$ws.Cells.Item(8, 17).Value = 0.00018288720041739672
$ws.Cells.Item(8, 18).Value = 0.0002714482848128302
$ws.Cells.Item(8, 19).Value = 0.0007045220174611162
$ws.Cells.Item(8, 20).Value = 0.00044697491351118076
$ws.Cells.Item(9, 1).Value = 45
$ws.Cells.Item(9, 2).Value = 0.006658541223783071
$ws.Cells.Item(9, 3).Value = 0.009958801289692943
$ws.Cells.Item(9, 4).Value = 0.00030832511964708453
$ws.Cells.Item(9, 5).Value = 0.008754045328446947
$ws.Cells.Item(9, 6).Value = 0.0008373544045444183
$ws.Cells.Item(9, 7).Value = 0.007973328695848374
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 10).Value = 0.01833712346224521
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 0.0010787683690666662
$ws.Cells.Item(9, 13).Value = 0.0005319116460180789
$ws.Cells.Item(9, 14).Value = 0.006957074545830529
$ws.Cells.Item(9, 15).Value = 0.00016034720724248695
$ws.Cells.Item(9, 16).Value = 0
$ws.Cells.Item(9, 17).Value = 0.000053977750699421105
$ws.Cells.Item(9, 18).Value = 0.000023998589204111956
$ws.Cells.Item(9, 19).Value = 0.000189817216907323
$ws.Cells.Item(9, 20).Value = 0.00006611874672898039
$ws.Cells.Item(10, 1).Value = 60
$ws.Cells.Item(10, 2).Value = 0.0005418001376983831
$ws.Cells.Item(10, 3).Value = 0.001186350831259178
$ws.Cells.Item(10, 4).Value = 0.000005795308441376953
$ws.Cells.Item(10, 5).Value = 0.004596595866693529
$ws.Cells.Item(10, 6).Value = 0.0008391681393675319
$ws.Cells.Item(10, 7).Value = 0.0030541914864323817
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = 0.01419109645777641
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 0.00003475703906332022
$ws.Cells.Item(10, 13).Value = 0.000008917032363793767
$ws.Cells.Item(10, 14).Value = 0.00016377960342791375
$ws.Cells.Item(10, 15).Value = 0.000004445530325614594
$ws.Cells.Item(10, 16).Value = 0
$ws.Cells.Item(10, 17).Value = 0
$ws.Cells.Item(10, 18).Value = 0.0000009584899041075681
$ws.Cells.Item(10, 19).Value = 0.000005483938787799678
$ws.Cells.Item(10, 20).Value = 0.00005006902204736992
$ws.Cells.Item(11, 1).Value = 75
$ws.Cells.Item(11, 2).Value = 0.00036018395465300316
$ws.Cells.Item(11, 3).Value = 0.0008106632178617441
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = 0.0031608632750471724
$ws.Cells.Item(11, 6).Value = 0.0008523468050379857
$ws.Cells.Item(11, 7).Value = 0.003497899172377893
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 10).Value = 0.01416919573520996
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 0.000003533426693668802
$ws.Cells.Item(11, 13).Value = 0
$ws.Cells.Item(11, 14).Value = 0.00002694857496008414
$ws.Cells.Item(11, 15).Value = 0
$ws.Cells.Item(11, 16).Value = 0
$ws.Cells.Item(11, 17).Value = 0
$ws.Cells.Item(11, 18).Value = 0
$ws.Cells.Item(11, 19).Value = 0
$ws.Cells.Item(11, 20).Value = 0.000029073021523452665
$ws.Cells.Item(12, 1).Value = 100
$ws.Cells.Item(12, 2).Value = 0.00017637248782687332
$ws.Cells.Item(12, 3).Value = 0.00039786479109630527
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = 0.0014472191713162117
$ws.Cells.Item(12, 6).Value = 0.0008559238101290488
$ws.Cells.Item(12, 7).Value = 0.0038523539927161447
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 0.014271504244954195
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 0.0000006205345524398561
$ws.Cells.Item(12, 13).Value = 0
$ws.Cells.Item(12, 14).Value = 0.000003460908395288108
$ws.Cells.Item(12, 15).Value = 0
$ws.Cells.Item(12, 16).Value = 0
$ws.Cells.Item(12, 17).Value = 0
$ws.Cells.Item(12, 18).Value = 0
$ws.Cells.Item(12, 19).Value = 0
$ws.Cells.Item(12, 20).Value = 0.000013431281438449574

